# Atualização 08/07 - Correção de bug: Uma nova planilha era criada sempre
# que o programa era aberto novamente.
#
# Instead of creating a brand-new workbook/sheet every time the clock-punch
# app is opened, it now reuses the existing worksheet: the most recent
# punch row is refreshed with the latest times and further punches are
# appended as new rows below it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the existing row 2 (ENTRADA/INTERVALO/RETORNO INTERVALO/SAÍDA)
# with the new punch times recorded for this session.
$ws.Cells.Item(2, 2).Value = "12:31:21"
$ws.Cells.Item(2, 3).Value = "12:31:22"
$ws.Cells.Item(2, 4).Value = "12:31:23"
$ws.Cells.Item(2, 5).Value = "12:31:24"

# Row 3: another punch cycle, same date, logged a bit later. Copying the
# date cell (instead of typing "08/07/2023" into .Value) keeps it a plain
# text value — typing it would make Excel auto-detect a date and assign a
# new number-format style to the cell, which the original file doesn't have.
$ws.Cells.Item(2, 1).Copy($ws.Cells.Item(3, 1))
$ws.Cells.Item(3, 2).Value = "12:31:43"
$ws.Cells.Item(3, 3).Value = "12:31:44"
$ws.Cells.Item(3, 4).Value = "12:31:45"
$ws.Cells.Item(3, 5).Value = "12:31:46"

# Row 4: yet another punch cycle logged later still.
$ws.Cells.Item(2, 1).Copy($ws.Cells.Item(4, 1))
$ws.Cells.Item(4, 2).Value = "12:32:03"
$ws.Cells.Item(4, 3).Value = "12:32:04"
$ws.Cells.Item(4, 4).Value = "12:32:05"
$ws.Cells.Item(4, 5).Value = "12:32:06"
